# Update "想去人数" (F) and "最低票价" (G) figures for several events on both
# the "展览" and "全部类型" worksheets, matching the refreshed scrape output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F5").Value = 6737
$ws1.Range("F6").Value = 87
$ws1.Range("F8").Value = 141
$ws1.Range("F9").Value = 6262
$ws1.Range("F10").Value = 49
$ws1.Range("F12").Value = 1261
$ws1.Range("G12").Value = 19.9
$ws1.Range("F21").Value = 4597
$ws1.Range("F23").Value = 41
$ws1.Range("F24").Value = 52
$ws1.Range("F25").Value = 194
$ws1.Range("F26").Value = 71

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F5").Value = 6737
$ws4.Range("F6").Value = 87
$ws4.Range("F8").Value = 141
$ws4.Range("F9").Value = 6262
$ws4.Range("F10").Value = 49
$ws4.Range("F12").Value = 1261
$ws4.Range("G12").Value = 19.9
$ws4.Range("F21").Value = 4597
$ws4.Range("F24").Value = 41
$ws4.Range("F25").Value = 52
$ws4.Range("F26").Value = 194
$ws4.Range("F27").Value = 71
